# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '44.020.57'
Set-TextCell 'E2' '  +0.05%  '
Set-TextCell 'D3' '2.238.46'
Set-TextCell 'E3' '  -0.24%  '
Set-TextCell 'E4' '  +0.21%  '
Set-TextCell 'D5' '305.46'
Set-TextCell 'E5' '  -3.89%  '
Set-TextCell 'D6' '94.40'
Set-TextCell 'E6' '  -6.14%  '
Set-TextCell 'E7' '  -0.48%  '
Set-TextCell 'E8' '  +0.25%  '
Set-TextCell 'D9' '0.523'
Set-TextCell 'E9' '  -3.59%  '
Set-TextCell 'D10' '34.74'
Set-TextCell 'E10' '  -5.50%  '
Set-TextCell 'D11' '0.0808'
Set-TextCell 'E11' '  -2.17%  '
Set-TextCell 'D12' '7.20'
Set-TextCell 'E12' '  -4.09%  '
Set-TextCell 'E13' '  -0.71%  '
Set-TextCell 'D14' '2.575.66'
Set-TextCell 'E14' '  -0.35%  '
Set-TextCell 'D15' '2.241.35'
Set-TextCell 'E15' '  -0.27%  '
Set-TextCell 'D16' '0.821'
Set-TextCell 'E16' '  -2.83%  '
Set-TextCell 'D17' '13.49'
Set-TextCell 'E17' '  -4.41%  '
Set-TextCell 'D18' '43.872.25'
Set-TextCell 'E18' '  -0.06%  '
Set-TextCell 'D19' '0.0₃0962'
Set-TextCell 'E19' '  -1.01%  '
Set-TextCell 'D20' '12.11'
Set-TextCell 'E20' '  -8.90%  '
Set-TextCell 'D21' '6.30'
Set-TextCell 'E21' '  -1.91%  '
Set-TextCell 'D22' '65.55'
Set-TextCell 'E22' '  +0.16%  '
Set-TextCell 'D23' '236.68'
Set-TextCell 'E23' '  +0.87%  '
Set-TextCell 'D24' '2.91'
Set-TextCell 'E24' '  -5.65%  '
Set-TextCell 'D25' '1.97'
Set-TextCell 'E25' '  -4.21%  '
Set-TextCell 'E26' '  +0.14%  '
Set-TextCell 'B27' 'InjectiveProtocol'
Set-TextCell 'C27' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D27' '38.27'
Set-TextCell 'E27' '  +2.26%  '
Set-TextCell 'B28' 'Cosmos'
Set-TextCell 'C28' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D28' '9.87'
Set-TextCell 'E28' '  -5.52%  '
Set-TextCell 'B29' 'Toncoin'
Set-TextCell 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D29' '2.20'
Set-TextCell 'E29' '  -0.03%  '
Set-TextCell 'D30' '6.03'
Set-TextCell 'E30' '  -1.59%  '
Set-TextCell 'D31' '19.87'
Set-TextCell 'E31' '  -0.80%  '
Set-TextCell 'D32' '149.86'
Set-TextCell 'E32' '  -6.01%  '
Set-TextCell 'D33' '0.0798'
Set-TextCell 'E33' '  -5.51%  '
Set-TextCell 'E34' '  -2.67%  '
Set-TextCell 'D35' '3.16'
Set-TextCell 'E35' '  -0.49%  '
Set-TextCell 'E36' '  -3.37%  '
Set-TextCell 'E37' '  +0.93%  '
Set-TextCell 'D38' '1.77'
Set-TextCell 'E38' '  -8.81%  '
Set-TextCell 'D39' '14.99'
Set-TextCell 'E39' '  -7.25%  '
Set-TextCell 'B40' 'NEARProtocol'
Set-TextCell 'C40' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D40' '3.37'
Set-TextCell 'E40' '  -7.80%  '
Set-TextCell 'B41' 'RenderToken'
Set-TextCell 'C41' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D41' '3.82'
Set-TextCell 'E41' '  -7.09%  '
Set-TextCell 'E42' '  -5.94%  '
Set-TextCell 'E43' '  +0.42%  '
Set-TextCell 'D44' '1.738.14'
Set-TextCell 'E44' '  -0.30%  '
Set-TextCell 'D45' '84.83'
Set-TextCell 'E45' '  +3.99%  '
Set-TextCell 'D46' '0.187'
Set-TextCell 'E46' '  -4.91%  '
Set-TextCell 'D47' '99.87'
Set-TextCell 'E47' '  -2.16%  '
Set-TextCell 'D48' '4.94'
Set-TextCell 'E48' '  -3.85%  '
Set-TextCell 'D49' '8.08'
Set-TextCell 'E49' '  -1.57%  '
Set-TextCell 'D50' '68.80'
Set-TextCell 'E50' '  -7.30%  '
Set-TextCell 'D51' '53.71'
Set-TextCell 'E51' '  -6.35%  '
